$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

$ws.Range("D2").Value = "58.068.26"
$ws.Range("E2").Value = "  +0.19%  "

$ws.Range("D3").Value = "2.478.99"
$ws.Range("E3").Value = "  +0.57%  "

$ws.Range("E4").Value = "  +0.01%  "

Set-TextValue $ws.Range("D5") "518.55"
$ws.Range("E5").Value = "  +0.38%  "

Set-TextValue $ws.Range("D6") "131.57"
$ws.Range("E6").Value = "  +0.36%  "

Set-TextValue $ws.Range("D7") "0.997"
$ws.Range("E7").Value = "  -0.33%  "

Set-TextValue $ws.Range("D8") "0.554"
$ws.Range("E8").Value = "  -0.76%  "

$ws.Range("D9").Value = "2.510.27"
$ws.Range("E9").Value = "  +1.74%  "

Set-TextValue $ws.Range("D10") "0.0971"
$ws.Range("E10").Value = "  -1.85%  "

Set-TextValue $ws.Range("D12") "5.18"
$ws.Range("E12").Value = "  -2.46%  "

Set-TextValue $ws.Range("D13") "0.331"
$ws.Range("E13").Value = "  -2.63%  "

$ws.Range("D14").Value = "2.925.85"
$ws.Range("E14").Value = "  +0.70%  "

$ws.Range("D15").Value = "58.013.48"
$ws.Range("E15").Value = "  +0.16%  "

Set-TextValue $ws.Range("D16") "22.02"
$ws.Range("E16").Value = "  -1.18%  "

Set-TextValue $ws.Range("D17") "0.0000134"
$ws.Range("E17").Value = "  -1.09%  "

$ws.Range("D18").Value = "2.498.57"
$ws.Range("E18").Value = "  +0.94%  "

Set-TextValue $ws.Range("D19") "10.70"
$ws.Range("E19").Value = "  -0.14%  "

Set-TextValue $ws.Range("D20") "320.52"
$ws.Range("E20").Value = "  +0.26%  "

Set-TextValue $ws.Range("D21") "4.16"
$ws.Range("E21").Value = "  -0.18%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D22") "5.97"
$ws.Range("E22").Value = "  +4.57%  "

$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D23") "0.997"
$ws.Range("E23").Value = "  -0.34%  "

Set-TextValue $ws.Range("D24") "64.16"
$ws.Range("E24").Value = "  +0.00%  "

Set-TextValue $ws.Range("D25") "0.400"
$ws.Range("E25").Value = "  -2.10%  "

Set-TextValue $ws.Range("D26") "0.993"
$ws.Range("E26").Value = "  -0.64%  "

$ws.Range("E27").Value = "  +0.66%  "

Set-TextValue $ws.Range("D28") "7.30"
$ws.Range("E28").Value = "  -0.13%  "

$ws.Range("D29").Value = "0.0₃0745"
$ws.Range("E29").Value = "  -0.08%  "

Set-TextValue $ws.Range("D30") "167.72"
$ws.Range("E30").Value = "  +1.73%  "

Set-TextValue $ws.Range("D31") "1.70"
$ws.Range("E31").Value = "  +0.93%  "

Set-TextValue $ws.Range("D32") "1.17"
$ws.Range("E32").Value = "  +1.11%  "

Set-TextValue $ws.Range("D33") "6.22"
$ws.Range("E33").Value = "  -0.19%  "

$ws.Range("E34").Value = "  -0.11%  "

Set-TextValue $ws.Range("D35") "0.994"
$ws.Range("E35").Value = "  -0.60%  "

Set-TextValue $ws.Range("D36") "18.01"
$ws.Range("E36").Value = "  -0.05%  "

Set-TextValue $ws.Range("D37") "1.25"
$ws.Range("E37").Value = "  -3.12%  "

Set-TextValue $ws.Range("D38") "3.92"
$ws.Range("E38").Value = "  -1.11%  "

Set-TextValue $ws.Range("D39") "36.77"
$ws.Range("E39").Value = "  +0.61%  "

Set-TextValue $ws.Range("D40") "1.45"
$ws.Range("E40").Value = "  -1.10%  "

Set-TextValue $ws.Range("D41") "0.766"
$ws.Range("E41").Value = "  -2.80%  "

Set-TextValue $ws.Range("D42") "274.63"
$ws.Range("E42").Value = "  +1.25%  "

Set-TextValue $ws.Range("D43") "5.07"
$ws.Range("E43").Value = "  +2.26%  "

Set-TextValue $ws.Range("D44") "3.42"
$ws.Range("E44").Value = "  -0.60%  "

Set-TextValue $ws.Range("D45") "0.594"
$ws.Range("E45").Value = "  +0.59%  "

Set-TextValue $ws.Range("D46") "0.0919"
$ws.Range("E46").Value = "  +1.55%  "

Set-TextValue $ws.Range("D47") "120.72"
$ws.Range("E47").Value = "  -4.65%  "

Set-TextValue $ws.Range("D48") "0.0497"
$ws.Range("E48").Value = "  +2.22%  "

Set-TextValue $ws.Range("D49") "17.75"
$ws.Range("E49").Value = "  +0.11%  "

Set-TextValue $ws.Range("D50") "0.0212"
$ws.Range("E50").Value = "  +0.86%  "

Set-TextValue $ws.Range("D51") "16.78"
$ws.Range("E51").Value = "  -0.46%  "
